$d = $word.ActiveDocument

function FindRange([string]$searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        return $rng
    }
    return $null
}

# --- Change 1: merge "Submission Date" + bookmark + ": 09-10-2019 " into one run ---
# Re-setting identical text across the bookmarked region collapses the runs into a
# single run and drops the now-redundant "_GoBack" bookmark that sat between them.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Submission Date: 09-10-2019", $true, $false, $false, $false, $false, $true, 1, $false, "Submission Date: 09-10-2019", 2) | Out-Null

# --- Change 2: remove the extra blank paragraph before "1. State Handshaking..." ---
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "State Handshaking Theorem") {
        $targetIdx = $i
        break
    }
}
$d.Paragraphs.Item($targetIdx - 1).Range.Delete()

# --- Change 3: drop the "?" in "State Handshaking Theorem?" ---
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("State Handshaking Theorem? Verify", $true, $false, $false, $false, $false, $true, 1, $false, "State Handshaking Theorem Verify", 2) | Out-Null

# --- Change 5: split "degrees of every path of vertices vi," into three runs, turning "path" into "pair" ---
$pathRange = FindRange("degrees of every path of vertices vi,")
$pStart = $pathRange.Start
$splitA = $pStart + 19   # right after "degrees of every pa"
$splitB = $pStart + 21   # right after "...path" (i.e. after "th")

$tmpBm1 = $d.Bookmarks.Add("zzTempSplit1", $d.Range($splitA, $splitA))
$tmpBm2 = $d.Bookmarks.Add("zzTempSplit2", $d.Range($splitB, $splitB))

$midRange = $d.Range($splitA, $splitB)
$midRange.Text = "ir"

$d.Bookmarks.Item("zzTempSplit1").Delete()
$d.Bookmarks.Item("zzTempSplit2").Delete()

# --- Change 4: split "State and Prove Dirac's Theorem." and drop a "_GoBack" bookmark in between ---
$diracRange = FindRange("Dirac's Theorem.")
$splitPoint = $diracRange.Start
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint)) | Out-Null
